# Applies the data refresh from re-running the EDA report after the
# train/test split changed the row count (7044 -> 7043) and a few
# dependent summary statistics shifted slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Rows (D) and Values (E) columns for every data row (2..21)
for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 4).Value = 7043   # column D - Rows
    $ws.Cells.Item($row, 5).Value = 7043   # column E - Values
}

# Row 6 (tenure): mean shifted slightly
$ws.Cells.Item(6, 11).Value = 32.37        # column K - mean

# Row 19 (MonthlyCharges): mean shifted slightly
$ws.Cells.Item(19, 11).Value = 64.76000000000001  # column K - mean

# Row 20 (TotalCharges): unique count and summary stats shifted
$ws.Cells.Item(20, 6).Value = 6531                # column F - Unique
$ws.Cells.Item(20, 11).Value = 2281.92            # column K - mean
$ws.Cells.Item(20, 12).Value = 2265.27            # column L - sd
